$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 284
$ws.Cells.Item(3, 6).Value = 285
$ws.Cells.Item(4, 6).Value = 286
$ws.Cells.Item(5, 6).Value = 287
$ws.Cells.Item(5, 12).Value = "stimuli/img_ose78.png"
$ws.Cells.Item(5, 13).Value = 80.19444444444444
$ws.Cells.Item(5, 14).Value = 60.25
$ws.Cells.Item(5, 15).Value = 70.22222222222223
$ws.Cells.Item(5, 16).Value = 36
$ws.Cells.Item(5, 17).Value = 8
$ws.Cells.Item(5, 18).Value = 7
$ws.Cells.Item(5, 19).Value = 7
$ws.Cells.Item(5, 20).Value = 7
$ws.Cells.Item(5, 21).Value = 7
$ws.Cells.Item(5, 22).Value = 7
$ws.Cells.Item(6, 6).Value = 288
$ws.Cells.Item(6, 12).Value = "stimuli/img_5p2ql.png"
$ws.Cells.Item(6, 13).Value = 89.19565217391305
$ws.Cells.Item(6, 14).Value = 72.52173913043478
$ws.Cells.Item(6, 15).Value = 80.85869565217391
$ws.Cells.Item(6, 16).Value = 46
$ws.Cells.Item(6, 17).Value = 10
$ws.Cells.Item(6, 18).Value = 10
$ws.Cells.Item(6, 19).Value = 10
$ws.Cells.Item(6, 20).Value = 10
$ws.Cells.Item(6, 21).Value = 10
$ws.Cells.Item(6, 22).Value = 9
$ws.Cells.Item(7, 6).Value = 289
$ws.Cells.Item(8, 6).Value = 290
$ws.Cells.Item(8, 12).Value = "stimuli/img_gbypq.png"
$ws.Cells.Item(8, 13).Value = 76.27500000000001
$ws.Cells.Item(8, 14).Value = 51.925
$ws.Cells.Item(8, 15).Value = 64.09999999999999
$ws.Cells.Item(8, 16).Value = 40
$ws.Cells.Item(8, 17).Value = 6
$ws.Cells.Item(8, 18).Value = 6
$ws.Cells.Item(8, 19).Value = 6
$ws.Cells.Item(8, 20).Value = 6
$ws.Cells.Item(8, 21).Value = 6
$ws.Cells.Item(8, 22).Value = 6
$ws.Cells.Item(9, 6).Value = 291
$ws.Cells.Item(9, 12).Value = "stimuli/img_aweye.png"
$ws.Cells.Item(9, 13).Value = 53.42105263157895
$ws.Cells.Item(9, 14).Value = 31.84210526315789
$ws.Cells.Item(9, 15).Value = 42.63157894736842
$ws.Cells.Item(9, 16).Value = 38
$ws.Cells.Item(9, 17).Value = 2
$ws.Cells.Item(9, 18).Value = 2
$ws.Cells.Item(9, 19).Value = 2
$ws.Cells.Item(9, 20).Value = 3
$ws.Cells.Item(9, 21).Value = 3
$ws.Cells.Item(9, 22).Value = 2
$ws.Cells.Item(10, 6).Value = 292
$ws.Cells.Item(11, 6).Value = 293
$ws.Cells.Item(12, 6).Value = 294
$ws.Cells.Item(12, 12).Value = "stimuli/img_bj2gr.png"
$ws.Cells.Item(12, 13).Value = 65.25
$ws.Cells.Item(12, 14).Value = 44.8
$ws.Cells.Item(12, 15).Value = 55.025
$ws.Cells.Item(12, 16).Value = 40
$ws.Cells.Item(12, 17).Value = 4
$ws.Cells.Item(12, 18).Value = 4
$ws.Cells.Item(12, 19).Value = 4
$ws.Cells.Item(12, 20).Value = 4
$ws.Cells.Item(12, 21).Value = 4
$ws.Cells.Item(12, 22).Value = 4
$ws.Cells.Item(13, 6).Value = 295
$ws.Cells.Item(14, 6).Value = 296
$ws.Cells.Item(15, 6).Value = 297
$ws.Cells.Item(16, 6).Value = 298
$ws.Cells.Item(16, 12).Value = "stimuli/img_okvvw.png"
$ws.Cells.Item(16, 13).Value = 50.58333333333334
$ws.Cells.Item(16, 14).Value = 32.11111111111111
$ws.Cells.Item(16, 15).Value = 41.34722222222223
$ws.Cells.Item(16, 16).Value = 36
$ws.Cells.Item(16, 17).Value = 2
$ws.Cells.Item(16, 18).Value = 2
$ws.Cells.Item(16, 19).Value = 2
$ws.Cells.Item(16, 20).Value = 2
$ws.Cells.Item(16, 21).Value = 2
$ws.Cells.Item(16, 22).Value = 3
$ws.Cells.Item(17, 6).Value = 299
$ws.Cells.Item(17, 12).Value = "stimuli/img_v8dra.png"
$ws.Cells.Item(17, 13).Value = 61.77272727272727
$ws.Cells.Item(17, 14).Value = 38.79545454545455
$ws.Cells.Item(17, 15).Value = 50.28409090909091
$ws.Cells.Item(17, 16).Value = 44
$ws.Cells.Item(17, 17).Value = 3
$ws.Cells.Item(17, 18).Value = 3
$ws.Cells.Item(17, 19).Value = 3
$ws.Cells.Item(17, 20).Value = 3
$ws.Cells.Item(17, 21).Value = 4
$ws.Cells.Item(18, 6).Value = 300
$ws.Cells.Item(19, 6).Value = 301
$ws.Cells.Item(19, 12).Value = "stimuli/img_t4hvr.png"
$ws.Cells.Item(19, 13).Value = 61.69230769230769
$ws.Cells.Item(19, 14).Value = 39.76923076923077
$ws.Cells.Item(19, 15).Value = 50.73076923076923
$ws.Cells.Item(19, 16).Value = 39
$ws.Cells.Item(19, 17).Value = 3
$ws.Cells.Item(19, 18).Value = 3
$ws.Cells.Item(19, 19).Value = 3
$ws.Cells.Item(19, 20).Value = 4
$ws.Cells.Item(19, 21).Value = 3
$ws.Cells.Item(19, 22).Value = 4
$ws.Cells.Item(20, 6).Value = 302
$ws.Cells.Item(21, 6).Value = 303
$ws.Cells.Item(21, 12).Value = "stimuli/img_kzg3h.png"
$ws.Cells.Item(21, 13).Value = 77.02777777777777
$ws.Cells.Item(21, 14).Value = 56.22222222222222
$ws.Cells.Item(21, 15).Value = 66.625
$ws.Cells.Item(21, 16).Value = 36
$ws.Cells.Item(21, 17).Value = 7
$ws.Cells.Item(21, 18).Value = 7
$ws.Cells.Item(21, 19).Value = 7
$ws.Cells.Item(21, 20).Value = 7
$ws.Cells.Item(21, 21).Value = 7
$ws.Cells.Item(21, 22).Value = 7
$ws.Cells.Item(22, 6).Value = 304
$ws.Cells.Item(23, 6).Value = 305
$ws.Cells.Item(24, 6).Value = 306
$ws.Cells.Item(24, 12).Value = "stimuli/img_2pk6v.png"
$ws.Cells.Item(24, 13).Value = 85.08108108108108
$ws.Cells.Item(24, 14).Value = 66.16216216216216
$ws.Cells.Item(24, 15).Value = 75.62162162162161
$ws.Cells.Item(24, 17).Value = 9
$ws.Cells.Item(24, 18).Value = 9
$ws.Cells.Item(24, 19).Value = 9
$ws.Cells.Item(24, 20).Value = 9
$ws.Cells.Item(24, 21).Value = 9
$ws.Cells.Item(24, 22).Value = 8
$ws.Cells.Item(25, 6).Value = 307
$ws.Cells.Item(26, 6).Value = 308
$ws.Cells.Item(26, 12).Value = "stimuli/img_ic3os.png"
$ws.Cells.Item(26, 13).Value = 84.79069767441861
$ws.Cells.Item(26, 14).Value = 66.16279069767442
$ws.Cells.Item(26, 15).Value = 75.47674418604652
$ws.Cells.Item(26, 17).Value = 9
$ws.Cells.Item(26, 18).Value = 9
$ws.Cells.Item(26, 19).Value = 9
$ws.Cells.Item(26, 20).Value = 8
$ws.Cells.Item(26, 21).Value = 9
$ws.Cells.Item(26, 22).Value = 9
$ws.Cells.Item(27, 6).Value = 309
$ws.Cells.Item(27, 12).Value = "stimuli/img_z3yzz.png"
$ws.Cells.Item(27, 13).Value = 71.71052631578948
$ws.Cells.Item(27, 14).Value = 49.81578947368421
$ws.Cells.Item(27, 15).Value = 60.76315789473685
$ws.Cells.Item(27, 16).Value = 38
$ws.Cells.Item(27, 17).Value = 5
$ws.Cells.Item(27, 18).Value = 5
$ws.Cells.Item(27, 19).Value = 5
$ws.Cells.Item(27, 20).Value = 5
$ws.Cells.Item(27, 21).Value = 5
$ws.Cells.Item(27, 22).Value = 5
$ws.Cells.Item(28, 6).Value = 310
$ws.Cells.Item(29, 6).Value = 311
$ws.Cells.Item(30, 6).Value = 312
$ws.Cells.Item(30, 12).Value = "stimuli/img_2pnl2.png"
$ws.Cells.Item(30, 13).Value = 6.621621621621622
$ws.Cells.Item(30, 14).Value = 7.135135135135135
$ws.Cells.Item(30, 15).Value = 6.878378378378379
$ws.Cells.Item(30, 17).Value = 1
$ws.Cells.Item(30, 18).Value = 1
$ws.Cells.Item(30, 19).Value = 1
$ws.Cells.Item(30, 20).Value = 1
$ws.Cells.Item(30, 21).Value = 1
$ws.Cells.Item(30, 22).Value = 1
$ws.Cells.Item(31, 6).Value = 313
$ws.Cells.Item(32, 6).Value = 314
$ws.Cells.Item(32, 12).Value = "stimuli/img_anzgh.png"
$ws.Cells.Item(32, 13).Value = 75.10526315789474
$ws.Cells.Item(32, 14).Value = 55.76315789473684
$ws.Cells.Item(32, 15).Value = 65.43421052631579
$ws.Cells.Item(32, 16).Value = 38
$ws.Cells.Item(32, 17).Value = 6
$ws.Cells.Item(32, 18).Value = 6
$ws.Cells.Item(32, 19).Value = 6
$ws.Cells.Item(32, 20).Value = 6
$ws.Cells.Item(32, 21).Value = 6
$ws.Cells.Item(32, 22).Value = 6
$ws.Cells.Item(33, 6).Value = 315
$ws.Cells.Item(34, 6).Value = 316
$ws.Cells.Item(34, 12).Value = "stimuli/img_jivhq.png"
$ws.Cells.Item(34, 13).Value = 37
$ws.Cells.Item(34, 14).Value = 22.26530612244898
$ws.Cells.Item(34, 15).Value = 29.63265306122449
$ws.Cells.Item(34, 16).Value = 49
$ws.Cells.Item(34, 17).Value = 2
$ws.Cells.Item(34, 18).Value = 2
$ws.Cells.Item(34, 19).Value = 2
$ws.Cells.Item(34, 20).Value = 2
$ws.Cells.Item(34, 21).Value = 2
$ws.Cells.Item(34, 22).Value = 2
$ws.Cells.Item(35, 6).Value = 317
$ws.Cells.Item(35, 12).Value = "stimuli/img_fqgem.png"
$ws.Cells.Item(35, 13).Value = 80.75
$ws.Cells.Item(35, 14).Value = 61.475
$ws.Cells.Item(35, 15).Value = 71.1125
$ws.Cells.Item(35, 16).Value = 40
$ws.Cells.Item(35, 17).Value = 8
$ws.Cells.Item(35, 18).Value = 8
$ws.Cells.Item(35, 19).Value = 8
$ws.Cells.Item(35, 20).Value = 8
$ws.Cells.Item(35, 21).Value = 8
$ws.Cells.Item(35, 22).Value = 8
$ws.Cells.Item(36, 6).Value = 318
$ws.Cells.Item(37, 6).Value = 319
$ws.Cells.Item(38, 6).Value = 320
$ws.Cells.Item(39, 6).Value = 321
$ws.Cells.Item(39, 12).Value = "stimuli/img_yteqw.png"
$ws.Cells.Item(39, 13).Value = 66.83783783783784
$ws.Cells.Item(39, 14).Value = 43.78378378378378
$ws.Cells.Item(39, 15).Value = 55.31081081081081
$ws.Cells.Item(39, 16).Value = 37
$ws.Cells.Item(39, 17).Value = 4
$ws.Cells.Item(39, 18).Value = 4
$ws.Cells.Item(39, 19).Value = 4
$ws.Cells.Item(39, 20).Value = 5
$ws.Cells.Item(39, 21).Value = 4
$ws.Cells.Item(39, 22).Value = 4
$ws.Cells.Item(40, 6).Value = 322
$ws.Cells.Item(40, 12).Value = "stimuli/img_3bxjb.png"
$ws.Cells.Item(40, 13).Value = 87.28571428571429
$ws.Cells.Item(40, 14).Value = 72.65714285714286
$ws.Cells.Item(40, 15).Value = 79.97142857142858
$ws.Cells.Item(40, 16).Value = 35
$ws.Cells.Item(40, 17).Value = 10
$ws.Cells.Item(40, 18).Value = 10
$ws.Cells.Item(40, 19).Value = 10
$ws.Cells.Item(40, 20).Value = 9
$ws.Cells.Item(40, 21).Value = 9
$ws.Cells.Item(40, 22).Value = 10
$ws.Cells.Item(41, 6).Value = 323
$ws.Cells.Item(42, 6).Value = 324
$ws.Cells.Item(42, 12).Value = "stimuli/img_cgdyc.png"
$ws.Cells.Item(42, 13).Value = 32.93023255813954
$ws.Cells.Item(42, 14).Value = 14.04651162790698
$ws.Cells.Item(42, 15).Value = 23.48837209302326
$ws.Cells.Item(42, 16).Value = 43
$ws.Cells.Item(42, 17).Value = 1
$ws.Cells.Item(42, 18).Value = 1
$ws.Cells.Item(42, 19).Value = 1
$ws.Cells.Item(42, 20).Value = 1
$ws.Cells.Item(42, 21).Value = 1
$ws.Cells.Item(42, 22).Value = 1
